$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell value: A19 goes from FALSE to TRUE
$ws.Range("A19").Value = $true

# Update the view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("A12").Select()
